$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values in columns D/E are plain text (prices/percentages sourced from a
# scrape) but many look numeric to Excel (e.g. "1.001", "0.07870"). Assigning
# them directly through .Value would let Excel auto-convert them to numbers
# and silently drop formatting such as trailing zeros. Prefixing with a
# leading apostrophe forces text entry, exactly like typing them in the UI,
# while leaving the cell's NumberFormat as "General" (unchanged).

$ws.Range("D2").Value = "'27.980.85"
$ws.Range("E2").Value = "  -2.23%  "
$ws.Range("D3").Value = "'1.829.66"
$ws.Range("E3").Value = "  -1.11%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'324.13"
$ws.Range("E5").Value = "  -3.19%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").Value = "'0.4648"
$ws.Range("E7").Value = "  -0.30%  "
$ws.Range("D8").Value = "'0.3871"
$ws.Range("E8").Value = "  -1.23%  "
$ws.Range("D9").Value = "'0.07870"
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("D10").Value = "'0.9588"
$ws.Range("E10").Value = "  -2.73%  "
$ws.Range("D11").Value = "'21.83"
$ws.Range("E11").Value = "  -1.46%  "
$ws.Range("D12").Value = "'1.792.59"
$ws.Range("E12").Value = "  -7.21%  "
$ws.Range("D13").Value = "'5.665"
$ws.Range("E13").Value = "  -3.28%  "
$ws.Range("D14").Value = "'6.897"
$ws.Range("E14").Value = "  -1.80%  "
$ws.Range("D15").Value = "'0.06800"
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("D16").Value = "'87.25"
$ws.Range("E16").Value = "  -0.46%  "
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").Value = "'0.000009917"
$ws.Range("E18").Value = "  -1.68%  "
$ws.Range("D19").Value = "'16.57"
$ws.Range("E19").Value = "  -2.71%  "
$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = "  -0.02%  "
$ws.Range("D21").Value = "'27.981.36"
$ws.Range("E21").Value = "  -2.28%  "
$ws.Range("D22").Value = "'5.316"
$ws.Range("E22").Value = "  -1.62%  "
$ws.Range("D23").Value = "'10.97"
$ws.Range("E23").Value = "  -2.23%  "
$ws.Range("D24").Value = "'2.086"
$ws.Range("E24").Value = "  -2.08%  "
$ws.Range("D25").Value = "'2.036.36"
$ws.Range("E25").Value = "  -7.13%  "
$ws.Range("D26").Value = "'153.67"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("D27").Value = "'19.14"
$ws.Range("E27").Value = "  -1.43%  "
$ws.Range("D28").Value = "'5.719"
$ws.Range("E28").Value = "  -6.48%  "
$ws.Range("D29").Value = "'1.966"
$ws.Range("E29").Value = "  -2.73%  "
$ws.Range("D30").Value = "'117.40"
$ws.Range("E30").Value = "  -0.11%  "
$ws.Range("D31").Value = "'0.09255"
$ws.Range("E31").Value = "  -1.81%  "
$ws.Range("D32").Value = "'0.9308"
$ws.Range("E32").Value = "  -4.58%  "
$ws.Range("D33").Value = "'5.283"
$ws.Range("E33").Value = "  -1.72%  "
$ws.Range("E34").Value = "  -2.33%  "
$ws.Range("D35").Value = "'3.292"
$ws.Range("E35").Value = "  -5.87%  "
$ws.Range("D36").Value = "'0.05854"
$ws.Range("E36").Value = "  -4.41%  "
$ws.Range("D37").Value = "'0.02138"
$ws.Range("E37").Value = "  -2.72%  "
$ws.Range("D38").Value = "'1.145"
$ws.Range("E38").Value = "  -1.36%  "
$ws.Range("D39").Value = "'7.766"
$ws.Range("E39").Value = "  +2.10%  "
$ws.Range("D40").Value = "'0.5575"
$ws.Range("E40").Value = "  -2.29%  "
$ws.Range("D41").Value = "'9.860"
$ws.Range("E41").Value = "  -2.31%  "
$ws.Range("D42").Value = "'0.1759"
$ws.Range("E42").Value = "  -1.98%  "
$ws.Range("D43").Value = "'11.63"
$ws.Range("E43").Value = "  -1.48%  "
$ws.Range("D44").Value = "'0.5259"
$ws.Range("E44").Value = "  -2.51%  "
$ws.Range("D45").Value = "'0.07003"
$ws.Range("E45").Value = "  -2.36%  "
$ws.Range("D46").Value = "'2.132"
$ws.Range("E46").Value = "  -10.41%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'1.825"
$ws.Range("E47").Value = "  -4.37%  "
$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'113.05"
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "'1.103"
$ws.Range("E49").Value = "  -11.99%  "
$ws.Range("D50").Value = "'1.001"
$ws.Range("E50").Value = "  +0.01%  "
$ws.Range("D51").Value = "'2.321"
$ws.Range("E51").Value = "  +0.31%  "
